$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.089.93"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "2.383.44"
$ws.Range("E3").Value = "  +6.85%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.96"
$ws.Range("E5").Value = "  +9.82%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.55"
$ws.Range("E6").Value = "  -6.73%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +3.35%  "

# Row 8
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("E9").Value = "  +7.39%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.75"
$ws.Range("E10").Value = "  -5.08%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("E11").Value = "  +2.25%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.58"
$ws.Range("E12").Value = "  -1.27%  "

# Row 13
$ws.Range("E13").Value = "  -4.58%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.21"
$ws.Range("E14").Value = "  +14.50%  "

# Row 16
$ws.Range("D16").Value = "2.745.95"
$ws.Range("E16").Value = "  +7.01%  "

# Row 17
$ws.Range("D17").Value = "2.405.49"
$ws.Range("E17").Value = "  +7.30%  "

# Row 18
$ws.Range("D18").Value = "43.193.64"
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.70"
$ws.Range("E19").Value = "  +6.31%  "

# Row 20
$ws.Range("E20").Value = "  +2.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "76.70"
$ws.Range("E21").Value = "  +4.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "274.52"
$ws.Range("E22").Value = "  +16.28%  "

# Row 23
$ws.Range("E23").Value = "  +2.06%  "

# Row 24
$ws.Range("E24").Value = "  -0.17%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.74"
$ws.Range("E25").Value = "  +8.74%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.74"
$ws.Range("E26").Value = "  +2.18%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28
$ws.Range("E28").Value = "  +6.90%  "

# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "176.64"
$ws.Range("E29").Value = "  +0.72%  "

# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -0.22%  "

# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.88"
$ws.Range("E31").Value = "  +0.59%  "

# Row 32
$ws.Range("E32").Value = "  +2.42%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0925"
$ws.Range("E33").Value = "  +4.29%  "

# Row 34
$ws.Range("E34").Value = "  +3.81%  "

# Row 35
$ws.Range("E35").Value = "  +5.27%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.86"
$ws.Range("E36").Value = "  -3.52%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.11"
$ws.Range("E37").Value = "  -2.71%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0366"
$ws.Range("E38").Value = "  -2.77%  "

# Row 39
$ws.Range("E39").Value = "  +2.77%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.83"
$ws.Range("E40").Value = "  +17.87%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.58"
$ws.Range("E41").Value = "  +20.77%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.232"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "124.34"
$ws.Range("E43").Value = "  +21.84%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.14"
$ws.Range("E44").Value = "  -4.01%  "

# Row 45
$ws.Range("E45").Value = "  +0.20%  "

# Row 46
$ws.Range("B46").Value = "BitcoinSV"
$ws.Range("C46").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "93.04"
$ws.Range("E46").Value = "  +67.68%  "

# Row 47
$ws.Range("B47").Value = "Celestia"
$ws.Range("C47").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.49"
$ws.Range("E47").Value = "  +0.53%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.52"
$ws.Range("E48").Value = "  +12.64%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.61"
$ws.Range("E49").Value = "  +3.93%  "

# Row 50
$ws.Range("E50").Value = "  +1.65%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.103"
$ws.Range("E51").Value = "  +5.42%  "
